$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-6 with new values, in the same order the original
# author entered them so the shared-strings table ends up in the same order.
$ws.Range("C2").Value = "11090"

$ws.Range("A3").Value = "John"
$ws.Range("B3").Value = "Doe"
$ws.Range("C3").Value = "11323"

$ws.Range("A4").Value = "Mikel"
$ws.Range("B4").Value = "Andjelo"
$ws.Range("C4").Value = "123123"

$ws.Range("A5").Value = "Zorz"
$ws.Range("B5").Value = "Zorz"

$ws.Range("A6").Value = "Sam"
$ws.Range("B6").Value = "Smith"
$ws.Range("C6").Value = "123122"

$ws.Range("C5").Value = "223412"

# Update selection to C5
$ws.Range("C5").Select()

# Update window position
$excel.ActiveWindow.Left = 20025
$excel.ActiveWindow.Top = 855
